$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.242.01"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.819.45"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'313.04"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4455"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "'0.3761"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").Value = "'0.07388"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'0.8791"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "'20.83"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.819.75"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'6.699"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "'5.411"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "'0.07116"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'0.000008805"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'15.01"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "27.260.73"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "'5.348"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'1.961"
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").Value = "'151.02"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "'2.298"
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("D27").Value = "'18.55"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'5.337"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "'117.12"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "'0.08861"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'0.7836"
$ws.Range("E31").Value = "  +4.06%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "'4.553"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("D34").Value = "'2.910"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'1.108"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "'0.01966"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.05258"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'7.292"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").Value = "'2.865"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'0.1703"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'2.282"
$ws.Range("E43").Value = "  +15.65%  "
$ws.Range("D44").Value = "'8.592"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").Value = "'0.5028"
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("D46").Value = "'10.56"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "'104.84"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'0.06381"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'66.02"
$ws.Range("E51").Value = "  +4.63%  "
